$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 89.12406158447266
$ws.Range("C2").Value = 7.183907985687256
$ws.Range("D2").Value = 47.98871994018555
$ws.Range("E2").Value = 57.85714340209961
